$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = 'LJ Speech'

$data = @(
    ,@(2, '<and eventually lodged in bethlehem hospital, where she died after forty years# detention, at the advanced age of one hundred.>', '<and he stention beforty years detention where she untrotion wheres detention where she untrod, at the efforty years detention thill lidd.>', 58.28624143467437)
    ,@(3, '<three months prior to his regularly scheduled separation date, ostensibly to care for his mother who had been injured in an accident at her work.>', '<three months prior who had been injured in anacts of leaod been injured in anaxient atter who had been injured in axy to care mother to his remonsibly scape.>', 58.9250161744424)
    ,@(4, '<in that period upwards of two hundred thousand writs>', '<in the nent those in then the nent thousand rich,>', 60.12513850200647)
    ,@(5, '<benavides stopped and waited in the truck until the gunman ran to the corner.>', '<benion ran to the corner until the gunman ran to the corner.>', 58.33690926765034)
    ,@(6, '<and to assume any disguise, if required to go into the country on the business of the association.>', '<anto the country cauntry cauntry canes of the usiness of the usiness of the usiness of the iron.>', 58.53982275362373)
    ,@(7, '<he denied that he had ever told frazier that he wanted a ride to irving to get curtain rods for an apartment.>', '<he diaver that he diartment to get curtain radpartment ride to get curving to get curving to get curtalen radpartment.>', 56.25656391218126)
    ,@(8, '<what the hell. i wonder what the hell is the uproar?>', '<one of the hallicr.>', 57.11770253634899)
    ,@(9, '<toward any official of the united states.>', '<to any of fishle of the enix.>', 58.86483345653676)
    ,@(10, '<in the marines, he evidenced a strong conviction as to the correctness of marxist doctrine,>', '<in the markedness of markedness of markedness of markedness he evidenced as to the ciren.>', 57.64059127529527)
    ,@(11, '<to the bitter destruction of all and everything. i am quite sure these people must hate not only the government>', '<to the begovernment every these people and every these people and every these people and every these people and every thing,>', 56.82235421269473)
    ,@(12, '<but they nevertheless frequently managed to secrete the means of making away with themselves, and accomplished their purpose.>', '<the less frequently making away with thems of making away with thems of making away with they with the making away with themselves.>', 56.76114451048021)
    ,@(13, '<in the case of convicted murderers only was prompt punishment inflicted,>', '<in the case of convicted murderers ow legted murderes ow legted.>', 59.04875597847063)
    ,@(14, '<#three# transports awaiting removal to the colonies, and #four# capital convicts, condemned to death and awaiting execution.>', '<fore, cat the no way and fore, cat the now way capital condemnt to deathen oway and forts condemnt to deaths, condemnt to deaths,>', 61.12224276971887)
    ,@(15, '<a great change must have occurred in oswald#s thinking to induce him to return to the united states.>', '<a great change must of a gretern to return to return to return to return to return to return to the uniting to states.>', 58.7625710417103)
    ,@(16, '<he was himself a prominent member of the low church, of austere piety, active in all good works.>', '<he was him in all good work stearch, of all good work stearch, activen member of all good works.>', 59.1798104866472)
    ,@(17, '<may not endorse him in #sixtyfour.>', '<main not idendor simmenciments iminsixty fore,>', 56.53264904739368)
    ,@(18, '<i was startled by the sharp report or explosion,>', '<our explosion,>', 57.4805352105398)
    ,@(19, '<he was coachman to a gentleman at roehampton. one day he went into a pawnbroker#s at wandsworth, and bought a pair of breeches on credit.>', '<he went day he went a bounds on credica brokers at reaches, undalton, at reaches, undalkers at reaches, on day he went to a paparive broaches at one dabbt.>', 58.4695198142644)
    ,@(20, '<and, if the council is used, arrangements should be made for the attendance of the secretary of the treasury>', '<and insould secretary of the triagements of the secretary of the triangements of the triangements of the trise used,>', 55.86536858936071)
    ,@(21, '<some of these decorations, the most valuable objects found in the ruins of the great city, still remain in their places on the walls#>', '<some of these deprations of these deprations of these deprations of these deprations of these deprations of these deprations of of thes.>', 60.47201498396409)
    ,@(22, '<oswald was five foot nine inches, slender and twentyfour years old.>', '<oswald was five four years alled was five fort ninetifors alled.>', 62.96605638364116)
    ,@(23, '<the service prefers to have two agents perform advance preparations.>', '<the service profers to have two advance prefersed ans preformance.>', 59.24487650910098)
    ,@(24, '<however, it now laid down in plain language and with precise details the requirements of a good jail system.>', '<house the requate jails details, the requirements of a good jails, holever admond to jails details details down.>', 58.19778193453484)
    ,@(25, '<the game, however, was nearly up. miss turner was met by her uncle on calais pier as she was walking with wakefield.>', '<the game with was a should callay per, a should callay per, a should callay perendar was nearly up, the ame how have and a opeled.>', 59.58386963479977)
    ,@(26, '<when he had taken the oath he was dismissed with the words, #get thee hence, wretch##>', '<when he own he had jed surred smisted with the words, get the words, get the words retch.>', 58.57370370739729)
    ,@(27, '<as was a changeofaddress card dated may twelve, nineteen sixtythree,>', '<as was at was at was at was at was at was at was at wellvelf night,>', 56.83706242043705)
    ,@(28, '<five point six, and six point five seconds.>', '<five poinc, and six poinc, and six boinds.>', 56.40332352957572)
    ,@(29, '<oswald#s revolver>', '<oswald#s rover>', 57.92443866287031)
    ,@(30, '<wallace#s relations warned him against his liverpool friend,>', '<owald assist has reporned. him igaganied.>', 56.08719090639899)
    ,@(31, '<this committee animadverted strongly upon the system in force at the metropolitan jails, and more especially upon the condition of newgate>', '<this commate, and more specially upon the strongly upon the strongly and were specially upon the condition of newgate, and more specially and were specite,>', 57.75621216624773)
    ,@(32, '<in response to inquiry by the commission regarding the instructions to agents in a motorcade>', '<in responstructions to inquigents and a motorcading the instructions to inquagents and emotorcade>', 58.85010745628207)
    ,@(33, '<to her husband#s membership card in the new orleans chapter.>', '<to orleans memober husband#s memobicard in the norleans memobins memopter.>', 58.20028447501357)
    ,@(34, '<who did not believe that he could reach, the president#s car at the speed it was then traveling.>', '<who dinot the speak who dinot that he could not be could not be could not be could reach,>', 61.69116487084581)
    ,@(35, '<his appearance on the scaffold was the signal for three prolonged shouts from an innumerable populace, #the brutal effusion of one common sentiment.#>', '<his opene common common commend to maral long for threen concomments entoman common common common common common common common concomment.>', 56.75843196206435)
    ,@(36, '<advised the bureau that oswald was unknown in such circles.>', '<advised that oswald was up that oswald was up yearo the areoths.>', 59.38592682089453)
    ,@(37, '<in the second simulated transcript which ended with the statement, quote, newspapers thank you, sir. you are a real patriot# end quote.>', '<in the seconded with the service than state and sire or real real real real real real real real real real real agent the statened the statened the statened thriellal.>', 58.71986451739383)
    ,@(38, '<quite half of the foregoing writs and arrests applied to sums under thirty pounds.>', '<a while to sipled you have of the forgoing rests applied thirty paunder thirty paunds.>', 57.73255860130507)
    ,@(39, '<to have obtained as much as forty thousand pounds by false and fraudulent pretenses.>', '<to havesent is forty thousand pounds by froms and falsem pounds and falsem pounds.>', 59.62767383130493)
    ,@(40, '<and told him a similar story to craig#s.>', '<and tol him a similar stract to crached or streate.>', 62.45816104884184)
    ,@(41, '<which grows as one speaks to him, end quote.>', '<which hrows as one speaks the him, end, quote.>', 59.20601155204875)
    ,@(42, '<they talked of the days when the convicts were hung up a dozen or more in a row#>', '<they told dust at adoe a more hone vortablo a dousen or more more more more more more more more.>', 55.52417503481846)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}

$ws.Range("A43:C52").ClearContents()